$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 38, shifting the existing rows 38-54 down to 39-55.
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with a new weekly price observation
# (same market/category series as the surrounding rows).
$ws.Range("A38").Value = 10
$ws.Range("B38").Value = "Vega Modelo de Temuco"
$ws.Range("C38").Value = "La Araucanía"
$ws.Range("D38").Value = 44466
$ws.Range("D38").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E38").Value = 9
$ws.Range("F38").Value = 100112035
$ws.Range("G38").Value = "Bruselas (repollito)"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 50
$ws.Range("K38").Value = 25000
$ws.Range("L38").Value = 25000
$ws.Range("M38").Value = 25000
$ws.Range("N38").Value = "`$/malla 10 kilos"
$ws.Range("O38").Value = "Provincia de Quillota"
$ws.Range("P38").Value = 2500
$ws.Range("Q38").Value = 10
$ws.Range("R38").Value = "Hortaliza"
